# Automatic update of files.
#
# Re-syncs the "Artfynd" sheet against an upstream export:
#   - rows 24/25/26 get their species-observation data rotated
#     (row24<-old row25, row25<-old row26, row26<-old row24)
#   - two additional observation rows (28, 29) are appended
#   - the sheet's used-range dimension grows from AY27 to AY29

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# ---- Row 24 (now holds what used to be row 25's data) ----
$ws.Range("A24").Value = 111986412
$ws.Range("Q24").Value = 396473.4754867578
$ws.Range("R24").Value = 6849402.350115799
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = "15:58"
$ws.Range("AB24").NumberFormat = "@"
$ws.Range("AB24").Value = "15:58"

# ---- Row 25 (now holds what used to be row 26's data) ----
$ws.Range("A25").Value = 111986331
$ws.Range("B25").Value = 90658
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 4361
$ws.Range("F25").Value = "Orange taggsvamp"
$ws.Range("G25").Value = "Hydnellum aurantiacum"
$ws.Range("H25").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q25").Value = 396461.0577280324
$ws.Range("R25").Value = 6849446.780411444
$ws.Range("Z25").NumberFormat = "@"
$ws.Range("Z25").Value = "15:50"
$ws.Range("AB25").NumberFormat = "@"
$ws.Range("AB25").Value = "15:50"

# ---- Row 26 (now holds what used to be row 24's data) ----
$ws.Range("A26").Value = 111986518
$ws.Range("B26").Value = 90678
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 4366
$ws.Range("F26").Value = "Skarp dropptaggsvamp"
$ws.Range("G26").Value = "Hydnellum peckii"
$ws.Range("H26").Value = "Banker"
$ws.Range("Q26").Value = 396445.8145670656
$ws.Range("R26").Value = 6849381.867442117
$ws.Range("Z26").NumberFormat = "@"
$ws.Range("Z26").Value = "16:06"
$ws.Range("AB26").NumberFormat = "@"
$ws.Range("AB26").Value = "16:06"

# ---- New row 28 ----
$ws.Range("A28").Value = 112060421
$ws.Range("B28").Value = 90660
$ws.Range("C28").Value = "Ovaliderad"
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 4362
$ws.Range("F28").Value = "Blå taggsvamp"
$ws.Range("G28").Value = "Hydnellum caeruleum"
$ws.Range("H28").Value = "(Hornem.) P.Karst."
$ws.Range("I28").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("P28").Value = "Mellandammen, Dlr"
$ws.Range("Q28").Value = 396290.5394533524
$ws.Range("R28").Value = 6849408.82232627
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = "Dalarna"
$ws.Range("U28").Value = "Älvdalen"
$ws.Range("V28").Value = "Dalarna"
$ws.Range("W28").Value = "Särna"
$ws.Range("Y28").NumberFormat = "@"
$ws.Range("Y28").Value = "2023-09-09"
$ws.Range("Z28").NumberFormat = "@"
$ws.Range("Z28").Value = "00:00"
$ws.Range("AA28").NumberFormat = "@"
$ws.Range("AA28").Value = "2023-09-09"
$ws.Range("AB28").NumberFormat = "@"
$ws.Range("AB28").Value = "00:00"
$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AF28").Value = ""
$ws.Range("AG28").Value = $false
$ws.Range("AT28").Value = ""
$ws.Range("AW28").Value = "Lars-Erik Nilsson"
$ws.Range("AX28").Value = "Lars-Erik Nilsson, Bo karlstens, Erik Danielsson, Göran Ehn, Håkan Thenander, Kajsa Larsson, Lisa Olson"
$ws.Range("AY28").Value = ""

# ---- New row 29 ----
$ws.Range("A29").Value = 112060422
$ws.Range("B29").Value = 90660
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 4362
$ws.Range("F29").Value = "Blå taggsvamp"
$ws.Range("G29").Value = "Hydnellum caeruleum"
$ws.Range("H29").Value = "(Hornem.) P.Karst."
$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = "3"
$ws.Range("J29").Value = "fruktkroppar"
$ws.Range("K29").Value = ""
$ws.Range("N29").Value = ""
$ws.Range("P29").Value = "Mellandammen, Dlr"
$ws.Range("Q29").Value = 396456.0580252151
$ws.Range("R29").Value = 6849454.036982723
$ws.Range("S29").Value = 5
$ws.Range("T29").Value = "Dalarna"
$ws.Range("U29").Value = "Älvdalen"
$ws.Range("V29").Value = "Dalarna"
$ws.Range("W29").Value = "Särna"
$ws.Range("Y29").NumberFormat = "@"
$ws.Range("Y29").Value = "2023-09-09"
$ws.Range("Z29").NumberFormat = "@"
$ws.Range("Z29").Value = "00:00"
$ws.Range("AA29").NumberFormat = "@"
$ws.Range("AA29").Value = "2023-09-09"
$ws.Range("AB29").NumberFormat = "@"
$ws.Range("AB29").Value = "00:00"
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AF29").Value = ""
$ws.Range("AG29").Value = $false
$ws.Range("AT29").Value = ""
$ws.Range("AW29").Value = "Lars-Erik Nilsson"
$ws.Range("AX29").Value = "Lars-Erik Nilsson, Bo karlstens, Erik Danielsson, Göran Ehn, Håkan Thenander, Kajsa Larsson, Lisa Olson"
$ws.Range("AY29").Value = ""
